$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.812.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.72%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.266.01'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '

# Row 6
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.25%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +4.58%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.53%  '

# Row 12
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.601.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.58%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.49%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.90%  '

# Row 16
$ws.Range("E16").Value = '  +3.68%  '

# Row 17
$ws.Range("E17").Value = '  +0.84%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.301.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.83%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.638.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.61%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.92%  '

# Row 22
$ws.Range("E22").Value = '  +2.56%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.62%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("E25").Value = '  +6.35%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '170.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.10%  '

# Row 29
$ws.Range("E29").Value = '  -1.69%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.62%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.38%  '

# Row 32
$ws.Range("E32").Value = '  +0.66%  '

# Row 33
$ws.Range("E33").Value = '  -0.10%  '

# Row 34
$ws.Range("E34").Value = '  +1.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.77'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.95%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0660'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.84%  '

# Row 37
$ws.Range("E37").Value = '  -2.63%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.93%  '

# Row 39
$ws.Range("E39").Value = '  -1.42%  '

# Row 40
$ws.Range("E40").Value = '  +4.01%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$ws.Range("E42").Value = '  +1.60%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000228'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.76%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0985'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '

# Row 45
$ws.Range("E45").Value = '  -7.06%  '

# Row 46
$ws.Range("E46").Value = '  -0.63%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.84%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.471.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.62%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.76%  '

# Row 50
$ws.Range("E50").Value = '  +0.66%  '

# Row 51
$ws.Range("E51").Value = '  +6.89%  '
